$d = $word.ActiveDocument

# The last paragraph currently ends with "TMDB API to get videos and movies and titles."
# Insert a brand new list paragraph right after it, inheriting the same
# paragraph/run formatting (ListParagraph style + numbering), then set its text.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "Used Swiper.js to create touch enabled and responsive sliders for the movies."
